# Applies the "Added Tantalum SMD Capacitors sizes: A B C D E V" change:
#  - fix typo "Capasitors_SMD" -> "Capacitors_SMD"
#  - add new worksheet "Tantalum_Capacitors_SMD" after the Capacitors sheet, with data
#  - make the new sheet the active tab/selection
#  - widen column J on the Resistors_SMD sheet slightly

$wb = $excel.ActiveWorkbook

$wsResistors  = $wb.Worksheets.Item(1)
$wsCapacitors = $wb.Worksheets.Item(2)

# --- fix the misspelled sheet name -------------------------------------
$wsCapacitors.Name = "Capacitors_SMD"

# --- widen column J (10) on the Resistors_SMD sheet ---------------------
$wsResistors.Columns.Item(10).ColumnWidth = 14.5

# --- add the new Tantalum_Capacitors_SMD worksheet at the end -----------
$wsTantalum = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCapacitors)
$wsTantalum.Name = "Tantalum_Capacitors_SMD"

# copy header-row formatting (fill/style) from the Capacitors sheet so the
# new sheet reuses the existing header style instead of creating a new one
$wsCapacitors.Range("A1:J1").Copy()
$wsTantalum.Range("A1:J1").PasteSpecial(-4122)

# header row (row 1) - reuse existing shared strings
$wsTantalum.Range("A1").Value = "Part Number"
$wsTantalum.Range("B1").Value = "Library Ref"
$wsTantalum.Range("C1").Value = "Footprint Ref"
$wsTantalum.Range("D1").Value = "Footprint Ref 2"
$wsTantalum.Range("E1").Value = "Footprint Ref 3"
$wsTantalum.Range("F1").Value = "Footprint Ref 4"
$wsTantalum.Range("G1").Value = "Footprint Ref 5"
$wsTantalum.Range("H1").Value = "Footprint Ref 6"
$wsTantalum.Range("I1").Value = "Library Path"
$wsTantalum.Range("J1").Value = "Footprint Path"

# data row (row 2) - set in the exact order so new shared strings are
# appended in the same sequence as the authored workbook
$wsTantalum.Range("A2").Value = "Tantalum Capacitor"
$wsTantalum.Range("B2").Value = "Polarized Capacitor"
$wsTantalum.Range("J2").Value = "CAP_SMD_TANTALUM.PcbLib"
$wsTantalum.Range("C2").Value = "TANT_CAP_SMD_A"
$wsTantalum.Range("D2").Value = "TANT_CAP_SMD_B"
$wsTantalum.Range("E2").Value = "TANT_CAP_SMD_C"
$wsTantalum.Range("F2").Value = "TANT_CAP_SMD_D"
$wsTantalum.Range("G2").Value = "TANT_CAP_SMD_E"
$wsTantalum.Range("H2").Value = "TANT_CAP_SMD_V"
$wsTantalum.Range("I2").Value = "Capacitors.SchLib"

# column widths (best-fit-like) for the new sheet
$wsTantalum.Columns.Item(1).ColumnWidth = 17.5
$wsTantalum.Columns.Item(2).ColumnWidth = 17.5
$wsTantalum.Columns.Item(3).ColumnWidth = 17
$wsTantalum.Columns.Item(4).ColumnWidth = 16.833333333333332
$wsTantalum.Columns.Item(5).ColumnWidth = 16.833333333333332
$wsTantalum.Columns.Item(6).ColumnWidth = 17
$wsTantalum.Columns.Item(7).ColumnWidth = 16.666666666666668
$wsTantalum.Columns.Item(8).ColumnWidth = 17
$wsTantalum.Columns.Item(9).ColumnWidth = 15.666666666666666
$wsTantalum.Columns.Item(10).ColumnWidth = 26.5

# --- selections / active sheet state ------------------------------------
# Capacitors_SMD is no longer the active tab; it keeps its frozen header
# pane but the cursor moves to H12.
[void]$wsCapacitors.Activate()
[void]$wsCapacitors.Range("H12").Select()

# Tantalum_Capacitors_SMD becomes the active/selected sheet with a frozen
# header row and the cursor on D6, matching a freshly duplicated sheet.
[void]$wsTantalum.Activate()
[void]$wsTantalum.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$wsTantalum.Range("D6").Select()
